$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2366412213740458
$ws.Range("C2").Value = 0.4732824427480916
$ws.Range("J2").Value = 0.04071246819338423
$ws.Range("P2").Value = 0.1501272264631043
$ws.Range("S2").Value = 0.09923664122137404
$ws.Range("B3").Value = 0.0106951871657754
$ws.Range("C3").Value = 0.0213903743315508
$ws.Range("J3").Value = 0.0748663101604278
$ws.Range("P3").Value = 0.6631016042780749
$ws.Range("S3").Value = 0.2299465240641711
$ws.Range("J4").Value = 0.04225352112676056
$ws.Range("P4").Value = 0.7183098591549296
$ws.Range("S4").Value = 0.2394366197183098
$ws.Range("B6").Value = 0.06372549019607843
$ws.Range("D6").Value = 0.009803921568627451
$ws.Range("E6").Value = 0.004901960784313725
$ws.Range("F6").Value = 0.04411764705882353
$ws.Range("J6").Value = 0.2107843137254902
$ws.Range("O6").Value = 0.03431372549019608
$ws.Range("Q6").Value = 0.1911764705882353
$ws.Range("R6").Value = 0.05392156862745098
$ws.Range("S6").Value = 0.3872549019607843
$ws.Range("B7").Value = 0.1330049261083744
$ws.Range("D7").Value = 0.03448275862068965
$ws.Range("F7").Value = 0.04926108374384237
$ws.Range("J7").Value = 0.1379310344827586
$ws.Range("O7").Value = 0.004926108374384237
$ws.Range("Q7").Value = 0.1822660098522167
$ws.Range("R7").Value = 0.03940886699507389
$ws.Range("S7").Value = 0.4187192118226601
$ws.Range("B8").Value = 0.099601593625498
$ws.Range("D8").Value = 0.0398406374501992
$ws.Range("F8").Value = 0.06772908366533864
$ws.Range("J8").Value = 0.08964143426294821
$ws.Range("O8").Value = 0.00597609561752988
$ws.Range("Q8").Value = 0.1852589641434263
$ws.Range("R8").Value = 0.1135458167330677
$ws.Range("S8").Value = 0.398406374501992
$ws.Range("B9").Value = 0.09774436090225563
$ws.Range("D9").Value = 0.02631578947368421
$ws.Range("E9").Value = 0.003759398496240601
$ws.Range("F9").Value = 0.06390977443609022
$ws.Range("J9").Value = 0.1240601503759398
$ws.Range("O9").Value = 0.007518796992481203
$ws.Range("Q9").Value = 0.1917293233082707
$ws.Range("R9").Value = 0.09022556390977443
$ws.Range("S9").Value = 0.3947368421052632
$ws.Range("B10").Value = 0.1346297681376215
$ws.Range("D10").Value = 0.02617801047120419
$ws.Range("F10").Value = 0.0643231114435303
$ws.Range("J10").Value = 0.1129394166043381
$ws.Range("O10").Value = 0.02692595362752431
$ws.Range("Q10").Value = 0.2101720269259536
$ws.Range("R10").Value = 0.07329842931937172
$ws.Range("S10").Value = 0.3515332834704563
$ws.Range("G11").Value = 0.1620795107033639
$ws.Range("J11").Value = 0.1009174311926606
$ws.Range("K11").Value = 0.2262996941896024
$ws.Range("L11").Value = 0.4892966360856269
$ws.Range("S11").Value = 0.02140672782874618
$ws.Range("G12").Value = 0.6932515337423313
$ws.Range("J12").Value = 0.2147239263803681
$ws.Range("L12").Value = 0.0245398773006135
$ws.Range("S12").Value = 0.06748466257668712
$ws.Range("G13").Value = 0.7017543859649122
$ws.Range("J13").Value = 0.2631578947368421
$ws.Range("S13").Value = 0.03508771929824561
$ws.Range("F15").Value = 0.004273504273504274
$ws.Range("H15").Value = 0.141025641025641
$ws.Range("I15").Value = 0.07692307692307693
$ws.Range("J15").Value = 0.3888888888888889
$ws.Range("K15").Value = 0.05555555555555555
$ws.Range("M15").Value = 0.01282051282051282
$ws.Range("O15").Value = 0.04700854700854701
$ws.Range("S15").Value = 0.2735042735042735
$ws.Range("F16").Value = 0.01276595744680851
$ws.Range("H16").Value = 0.174468085106383
$ws.Range("I16").Value = 0.05531914893617021
$ws.Range("J16").Value = 0.4042553191489361
$ws.Range("K16").Value = 0.1574468085106383
$ws.Range("M16").Value = 0.02553191489361702
$ws.Range("O16").Value = 0.05106382978723404
$ws.Range("S16").Value = 0.1191489361702128
$ws.Range("F17").Value = 0.01004016064257028
$ws.Range("H17").Value = 0.1686746987951807
$ws.Range("I17").Value = 0.1325301204819277
$ws.Range("J17").Value = 0.4116465863453815
$ws.Range("K17").Value = 0.08433734939759036
$ws.Range("M17").Value = 0.01405622489959839
$ws.Range("N17").Value = 0.002008032128514056
$ws.Range("O17").Value = 0.06224899598393574
$ws.Range("S17").Value = 0.1144578313253012
$ws.Range("F18").Value = 0.005025125628140704
$ws.Range("H18").Value = 0.1507537688442211
$ws.Range("I18").Value = 0.1206030150753769
$ws.Range("J18").Value = 0.3567839195979899
$ws.Range("K18").Value = 0.1306532663316583
$ws.Range("M18").Value = 0.03517587939698492
$ws.Range("O18").Value = 0.05025125628140704
$ws.Range("S18").Value = 0.1507537688442211
$ws.Range("F19").Value = 0.008708272859216255
$ws.Range("H19").Value = 0.227866473149492
$ws.Range("I19").Value = 0.1066763425253991
$ws.Range("J19").Value = 0.3505079825834543
$ws.Range("K19").Value = 0.09869375907111756
$ws.Range("M19").Value = 0.02685050798258345
$ws.Range("N19").Value = 0.000725689404934688
$ws.Range("O19").Value = 0.0660377358490566
$ws.Range("S19").Value = 0.113933236574746
